$wb = $excel.ActiveWorkbook

$confirmados = $wb.Worksheets.Item("Confirmados")
$mortes = $wb.Worksheets.Item("Mortes")

# New daily values (rows 48 and 49) for "Confirmados" (cumulative confirmed cases)
$confirmadosRow48 = @(72,48,193,1050,635,1582,579,349,209,344,121,100,750,217,85,676,816,41,2607,289,640,33,75,732,8419,42,23)
$confirmadosRow49 = @(77,48,230,1206,673,1676,614,383,229,398,123,101,806,246,101,738,960,44,2855,302,653,35,79,768,8755,44,25)

# New daily values (rows 48 and 49) for "Mortes" (cumulative deaths)
$mortesRow48 = @(2,3,3,53,21,67,14,9,10,21,3,2,17,10,11,26,72,7,155,13,15,2,3,21,560,4,0)
$mortesRow49 = @(2,3,5,62,21,74,14,9,14,24,3,2,20,13,13,30,85,7,170,15,16,2,3,24,588,4,0)

$confirmados.Range("A48").Value = "'2020-04-11"
$confirmados.Range("A49").Value = "'2020-04-12"
$mortes.Range("A48").Value = "'2020-04-11"
$mortes.Range("A49").Value = "'2020-04-12"

for ($i = 0; $i -lt $confirmadosRow48.Length; $i++) {
    $col = $i + 2
    $confirmados.Cells.Item(48, $col).Value = $confirmadosRow48[$i]
    $confirmados.Cells.Item(49, $col).Value = $confirmadosRow49[$i]
    $mortes.Cells.Item(48, $col).Value = $mortesRow48[$i]
    $mortes.Cells.Item(49, $col).Value = $mortesRow49[$i]
}
